$d = $word.ActiveDocument
$wns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# Process from the bottom of the document upward so earlier paragraph
# indices are never affected by edits made further down.

# Paragraph 22: "Alex Kokkosoulis wrote all of the UI. This includes all XAML files."
# + new trailing sentence about the uml and state diagram.
$inner22 = '<w:r><w:t>Alex Kokkosoulis wrote all of the UI. This includes all XAML files.</w:t></w:r>'
$inner22 += '<w:r><w:t xml:space="preserve"> The uml and state diagram was also completed.</w:t></w:r>'
$p22 = $d.Paragraphs(22)
$p22.Range.InsertXML('<w:p ' + $wns + '>' + $inner22 + '</w:p>')

# Paragraph 20: Nathaniel Barrett paragraph.
$inner20 = '<w:r><w:t xml:space="preserve">Nathaniel Barrett wrote the classes associated with the </w:t></w:r>'
$inner20 += '<w:r><w:t>OrderHistory and the Store selection.</w:t></w:r>'
$inner20 += '<w:r><w:t xml:space="preserve"> The uml diagram was also completed.</w:t></w:r>'
$p20 = $d.Paragraphs(20)
$p20.Range.InsertXML('<w:p ' + $wns + '>' + $inner20 + '</w:p>')

# Paragraph 18: Zachary Ellis paragraph (bookmark removed from here; sequence
# diagram sentence added at the end).
$inner18 = '<w:r><w:t>Zachary Ellis wrote the classes associated with the Settings menu</w:t></w:r>'
$inner18 += '<w:r><w:t xml:space="preserve"> and Unit-Testing</w:t></w:r>'
$inner18 += '<w:r><w:t>. This includes the serialization helper class, the Settings and ProfileSettings view models, and the UserSettings model class.</w:t></w:r>'
$inner18 += '<w:r><w:t xml:space="preserve"> </w:t></w:r>'
$inner18 += '<w:r><w:t>The sequence diagram was also completed.</w:t></w:r>'
$p18 = $d.Paragraphs(18)
$p18.Range.InsertXML('<w:p ' + $wns + '>' + $inner18 + '</w:p>')

# Paragraph 16: Jacob Damon paragraph with new trailing sentence.
$inner16 = '<w:r><w:t>Jacob Damon wrote the classes associated with the ordering system. This includes the product selection logic, the “Product” model classes, and the Cart/Checkout logic.</w:t></w:r>'
$inner16 += '<w:r><w:t xml:space="preserve"> The class, state, uml, and sequence diagrams were also completed.</w:t></w:r>'
$p16 = $d.Paragraphs(16)
$p16.Range.InsertXML('<w:p ' + $wns + '>' + $inner16 + '</w:p>')

# Paragraph 14: Jackson Dumas paragraph; the _GoBack bookmark moves here, and
# a new "SRRS" sentence is appended.
$inner14 = '<w:r><w:t>Jackson Dumas wrote the Web API, the SQL procedures, and the Firebase-Authentication implementation (used for logging in).</w:t></w:r>'
$inner14 += '<w:r><w:t xml:space="preserve"> All database functionality, platform-dependent Firebase classes, as well as the LoginViewModel class were written by Jackson.</w:t></w:r>'
$inner14 += '<w:r><w:t xml:space="preserve"> The SRRS was also completed.</w:t></w:r>'
$inner14 += '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>'
$p14 = $d.Paragraphs(14)
$p14.Range.InsertXML('<w:p ' + $wns + '>' + $inner14 + '</w:p>')

# Paragraph 10: Authors list entry "Alex Kokkosoulis" merged into a single run.
# Keep the existing list-paragraph formatting (ListParagraph style + bullet
# numbering) intact; only the runs themselves change.
$inner10 = '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr>'
$inner10 += '<w:r><w:t>Alex Kokkosoulis</w:t></w:r>'
$p10 = $d.Paragraphs(10)
$p10.Range.InsertXML('<w:p ' + $wns + '>' + $inner10 + '</w:p>')
